# Remove the now-unused "Sheet" reference row (row 16) from the
# optimization_parameters sheet, then restore the user's tab/selection
# state (network!J42 selected, threshold_b the active/displayed tab).

$wb = $excel.ActiveWorkbook

$wsOpt = $wb.Worksheets.Item("optimization_parameters")
$wsOpt.Rows.Item(16).Delete() | Out-Null
$wsOpt.Rows.Item(16).Select() | Out-Null

$wsNet = $wb.Worksheets.Item("network")
$wsNet.Activate() | Out-Null
$wsNet.Range("J42").Select() | Out-Null

$wsThresh = $wb.Worksheets.Item("threshold_b")
$wsThresh.Activate() | Out-Null
